# Fruta / hortaliza, semanal
#
# A new weekly price record (row) was added to the price log for
# "Terminal Hortofrutícola Agro Chillán - Pepino dulce". It belongs above
# the existing data (most-recent-first ordering), so we insert a new
# row at position 10 - pushing the former rows 10-20 down to 11-21 - and
# populate the new row 10 with the new reading.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 10; this shifts rows 10:20 down
# to 11:21 and carries the row's number formatting (e.g. the date style
# on column D) down with them, matching Excel's native Insert behaviour.
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row 10 with the new record.
$ws.Cells.Item(10, 1).Value = 7
$ws.Cells.Item(10, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(10, 3).Value = "Ñuble"
$ws.Cells.Item(10, 4).Value = 44775
$ws.Cells.Item(10, 5).Value = 16
$ws.Cells.Item(10, 6).Value = 100112043
$ws.Cells.Item(10, 7).Value = "Pepino dulce"
$ws.Cells.Item(10, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(10, 9).Value = "Primera"
$ws.Cells.Item(10, 10).Value = 100
$ws.Cells.Item(10, 11).Value = 17000
$ws.Cells.Item(10, 12).Value = 18000
$ws.Cells.Item(10, 13).Value = 17500
$ws.Cells.Item(10, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(10, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(10, 16).Value = 972
$ws.Cells.Item(10, 17).Value = 18
$ws.Cells.Item(10, 18).Value = "Hortaliza"
